# Append one new log row (row 25) to the Nalco run-log sheet, matching the
# "UPDATE Nalco PDF (2025-08-17 13:03:20 UTC)" commit: a SKIPPED run whose
# PDF didn't change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$lastRow = 24
$newRow = $lastRow + 1

# Same data pattern as the other "SKIPPED" rows: Run UTC, Run IST, Status,
# Message, Chosen URL, Saved PDF (blank), Rows Appended (0), Total Rows After (blank).
$ws.Cells.Item($newRow, 1).Value = "2025-08-17 13:03:19 UTC"
$ws.Cells.Item($newRow, 2).Value = "2025-08-17 18:33:19 IST"
$ws.Cells.Item($newRow, 3).Value = "SKIPPED"
$ws.Cells.Item($newRow, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($newRow, 5).Value = "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf"
$ws.Cells.Item($newRow, 6).Value = ""
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = ""

# Carry over the same cell styling (centered, no border) used by every
# other data row, by copying the formatting from the previous last row.
$ws.Range("A24:H24").Copy()
$ws.Range("A25:H25").PasteSpecial(-4122)
$excel.CutCopyMode = $false
